$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 20836216
$ws.Range("J113").Value = 62502500
$ws.Range("L113").Value = 62502500
$ws.Range("N113").Value = -62509008

$ws.Range("H116").Value = 6427.28
$ws.Range("I116").Value = 9498.462
$ws.Range("J116").Value = 3100.1667
$ws.Range("K116").Value = 9498.462
$ws.Range("L116").Value = 3100.1667
$ws.Range("M116").Value = -6056.462
$ws.Range("N116").Value = -9984.1667

$ws.Range("H129").Value = 955.75
$ws.Range("I129").Value = 732.25
$ws.Range("J129").Value = 982.0441
$ws.Range("K129").Value = 2196.75
$ws.Range("L129").Value = 2946.1323
$ws.Range("M129").Value = 2803.25
$ws.Range("N129").Value = -12946.1323

$ws.Range("H132").Value = 967.6769399999999
$ws.Range("I132").Value = 728.5685999999999
$ws.Range("K132").Value = 2185.7058
$ws.Range("M132").Value = 344.2942000000003

$ws.Range("H137").Value = 1032.2625
$ws.Range("I137").Value = 840.5857
$ws.Range("J137").Value = 2374
$ws.Range("K137").Value = 2521.7571
$ws.Range("L137").Value = 7122
$ws.Range("M137").Value = 28.24290000000019
$ws.Range("N137").Value = -12222

$ws.Range("H138").Value = 2850.2856
$ws.Range("I138").Value = 1223.6177
$ws.Range("J138").Value = 5364.227
$ws.Range("K138").Value = 3670.8531
$ws.Range("L138").Value = 16092.681
$ws.Range("M138").Value = 1469.1469
$ws.Range("N138").Value = -26372.681

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 59307
$ws.Range("I23").Value = 76672.664
$ws.Range("J23").Value = 7210
$ws.Range("K23").Value = 76672.664
$ws.Range("L23").Value = 7210
$ws.Range("M23").Value = -76413.664
$ws.Range("N23").Value = -7728

$ws.Range("H32").Value = 5512.2554
$ws.Range("I32").Value = 3789
$ws.Range("J32").Value = 14128.533
$ws.Range("K32").Value = 3789
$ws.Range("L32").Value = 14128.533
$ws.Range("M32").Value = -3502
$ws.Range("N32").Value = -14702.533

$ws.Range("H37").Value = 9197.125
$ws.Range("I37").Value = 7860
$ws.Range("J37").Value = 9999.4
$ws.Range("K37").Value = 7860
$ws.Range("L37").Value = 9999.4
$ws.Range("M37").Value = -7587
$ws.Range("N37").Value = -10545.4

$ws.Range("H61").Value = 4865.933
$ws.Range("I61").Value = 5371.92
$ws.Range("J61").Value = 2336
$ws.Range("K61").Value = 5371.92
$ws.Range("L61").Value = 2336
$ws.Range("M61").Value = -5159.92
$ws.Range("N61").Value = -2760

$ws.Range("H74").Value = 1444.6786
$ws.Range("I74").Value = 1421.5883
$ws.Range("J74").Value = 1480.3636
$ws.Range("K74").Value = 1421.5883
$ws.Range("L74").Value = 1480.3636
$ws.Range("M74").Value = -547.5882999999999
$ws.Range("N74").Value = -3228.3636

$ws.Range("H77").Value = 1444.6786
$ws.Range("I77").Value = 1421.5883
$ws.Range("J77").Value = 1480.3636
$ws.Range("K77").Value = 7107.941499999999
$ws.Range("L77").Value = 7401.817999999999
$ws.Range("M77").Value = -2739.941499999999
$ws.Range("N77").Value = -16137.818

$ws.Range("H102").Value = 9260684
$ws.Range("I102").Value = 9260684
$ws.Range("K102").Value = 9260684
$ws.Range("M102").Value = -9259062

$ws.Range("H110").Value = 1900.1666
$ws.Range("I110").Value = 880.2
$ws.Range("J110").Value = 7000
$ws.Range("K110").Value = 880.2
$ws.Range("L110").Value = 7000
$ws.Range("M110").Value = 1164.8
$ws.Range("N110").Value = -11090

$ws.Range("H122").Value = 1975419.1
$ws.Range("I122").Value = 3209431
$ws.Range("J122").Value = 999.8
$ws.Range("K122").Value = 9628293
$ws.Range("L122").Value = 2999.4
$ws.Range("M122").Value = -9625843
$ws.Range("N122").Value = -7899.4

$ws.Range("H132").Value = 1889203.6
$ws.Range("I132").Value = 1766.4615
$ws.Range("K132").Value = 5299.3845
$ws.Range("M132").Value = -2769.3845

$ws.Range("H136").Value = 4865.933
$ws.Range("I136").Value = 5371.92
$ws.Range("J136").Value = 2336
$ws.Range("K136").Value = 16115.76
$ws.Range("L136").Value = 7008
$ws.Range("M136").Value = -13565.76
$ws.Range("N136").Value = -12108

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1443.1666
$ws.Range("I94").Value = 474
$ws.Range("J94").Value = 2800
$ws.Range("K94").Value = 474
$ws.Range("L94").Value = 2800
$ws.Range("M94").Value = -23
$ws.Range("N94").Value = -3702

$ws.Range("H99").Value = 111112590
$ws.Range("I99").Value = 142858110
$ws.Range("K99").Value = 142858110
$ws.Range("M99").Value = -142856612

$ws.Range("H102").Value = 2406
$ws.Range("I102").Value = 2406
$ws.Range("K102").Value = 2406
$ws.Range("M102").Value = 839

$ws.Range("H107").Value = 142858660
$ws.Range("I107").Value = 166668100
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 166668100
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = -166666180
$ws.Range("N107").Value = -5840

$ws.Range("H134").Value = 3252.434
$ws.Range("I134").Value = 3290.848
$ws.Range("K134").Value = 9872.544
$ws.Range("M134").Value = -7337.544

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 384.92856
$ws.Range("I22").Value = 377.25
$ws.Range("J22").Value = 395.16666
$ws.Range("K22").Value = 377.25
$ws.Range("L22").Value = 395.16666
$ws.Range("M22").Value = -27.25
$ws.Range("N22").Value = -1095.16666

$ws.Range("H31").Value = 6748.361
$ws.Range("I31").Value = 1609.7949
$ws.Range("J31").Value = 15857.637
$ws.Range("K31").Value = 1609.7949
$ws.Range("L31").Value = 15857.637
$ws.Range("M31").Value = -1314.7949
$ws.Range("N31").Value = -16447.637

$ws.Range("H34").Value = 6748.361
$ws.Range("I34").Value = 1609.7949
$ws.Range("J34").Value = 15857.637
$ws.Range("K34").Value = 1609.7949
$ws.Range("L34").Value = 15857.637
$ws.Range("M34").Value = -1407.7949
$ws.Range("N34").Value = -16261.637

$ws.Range("H58").Value = 1024.3871
$ws.Range("I58").Value = 636.13635
$ws.Range("J58").Value = 1973.4445
$ws.Range("K58").Value = 636.13635
$ws.Range("L58").Value = 1973.4445
$ws.Range("M58").Value = -433.13635
$ws.Range("N58").Value = -2379.4445

$ws.Range("H99").Value = 11376147
$ws.Range("I99").Value = 13530.286
$ws.Range("K99").Value = 13530.286
$ws.Range("M99").Value = -12032.286

$ws.Range("H122").Value = 4065.25
$ws.Range("I122").Value = 3337.6667
$ws.Range("J122").Value = 4501.8
$ws.Range("K122").Value = 10013.0001
$ws.Range("L122").Value = 13505.4
$ws.Range("M122").Value = -7563.000100000001
$ws.Range("N122").Value = -18405.4

$ws.Range("H126").Value = 11376147
$ws.Range("I126").Value = 13530.286
$ws.Range("K126").Value = 40590.858
$ws.Range("M126").Value = -38120.858

$ws.Range("H132").Value = 1666.9608
$ws.Range("I132").Value = 1254.8788
$ws.Range("J132").Value = 2422.4443
$ws.Range("K132").Value = 3764.6364
$ws.Range("L132").Value = 7267.3329
$ws.Range("M132").Value = -1234.6364
$ws.Range("N132").Value = -12327.3329

$ws.Range("H134").Value = 1641.904
$ws.Range("I134").Value = 1816.7872
$ws.Range("J134").Value = 1325.7693
$ws.Range("K134").Value = 5450.3616
$ws.Range("L134").Value = 3977.3079
$ws.Range("M134").Value = -2915.3616
$ws.Range("N134").Value = -9047.3079

$ws.Range("H136").Value = 1024.3871
$ws.Range("I136").Value = 636.13635
$ws.Range("J136").Value = 1973.4445
$ws.Range("K136").Value = 1908.40905
$ws.Range("L136").Value = 5920.333500000001
$ws.Range("M136").Value = 641.59095
$ws.Range("N136").Value = -11020.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 357715.47
$ws.Range("I113").Value = 651.36365
$ws.Range("K113").Value = 1954.09095
$ws.Range("M113").Value = 215.90905

$ws.Range("H132").Value = 1482.125
$ws.Range("I132").Value = 693
$ws.Range("J132").Value = 3849.5
$ws.Range("K132").Value = 6237
$ws.Range("L132").Value = 34645.5
$ws.Range("M132").Value = -3707
$ws.Range("N132").Value = -39705.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1219.6052
$ws.Range("I102").Value = 902.4074000000001
$ws.Range("J102").Value = 1998.1818
$ws.Range("K102").Value = 902.4074000000001
$ws.Range("L102").Value = 1998.1818
$ws.Range("M102").Value = 719.5925999999999
$ws.Range("N102").Value = -5242.1818

$ws.Range("H132").Value = 1694.6923
$ws.Range("I132").Value = 1233.625
$ws.Range("J132").Value = 3231.5833
$ws.Range("K132").Value = 3700.875
$ws.Range("L132").Value = 9694.749899999999
$ws.Range("M132").Value = -1170.875
$ws.Range("N132").Value = -14754.7499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2168694.8
$ws.Range("I122").Value = 2861805
$ws.Range("J122").Value = 2725
$ws.Range("K122").Value = 8585415
$ws.Range("L122").Value = 8175
$ws.Range("M122").Value = -8582965
$ws.Range("N122").Value = -13075

$ws.Range("H132").Value = 9377683
$ws.Range("I132").Value = 12147644
$ws.Range("J132").Value = 2429.923
$ws.Range("K132").Value = 36442932
$ws.Range("L132").Value = 7289.768999999999
$ws.Range("M132").Value = -36440402
$ws.Range("N132").Value = -12349.769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1485.35
$ws.Range("I122").Value = 1420.8462
$ws.Range("J122").Value = 1605.1428
$ws.Range("K122").Value = 4262.5386
$ws.Range("L122").Value = 4815.428400000001
$ws.Range("M122").Value = -1812.5386
$ws.Range("N122").Value = -9715.428400000001
Write-Host "Applied all changes"
